$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.817.69"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.624.65"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.51"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.20"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0608"
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.856.81"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.616.97"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  -0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.96"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.844.26"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.90"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("E23").Value = "  -4.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.28"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.44"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.09"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.404.62"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("E38").Value = "  -0.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.553"
$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.845"
$ws.Range("E40").Value = "  -2.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -2.00%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.60"
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.765.68"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("E47").Value = "  -3.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.95"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  +1.19%  "
